# Auto-generated edit script applying F-column ("想去人数") value updates
# across the four worksheets, per the provided OOXML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value = 1981
$ws.Cells.Item(7, 6).Value = 1690
$ws.Cells.Item(8, 6).Value = 731
$ws.Cells.Item(13, 6).Value = 1764
$ws.Cells.Item(14, 6).Value = 1177
$ws.Cells.Item(15, 6).Value = 1697
$ws.Cells.Item(17, 6).Value = 25
$ws.Cells.Item(19, 6).Value = 1177
$ws.Cells.Item(20, 6).Value = 1783
$ws.Cells.Item(21, 6).Value = 1783
$ws.Cells.Item(22, 6).Value = 1185
$ws.Cells.Item(25, 6).Value = 1665
$ws.Cells.Item(27, 6).Value = 4497
$ws.Cells.Item(29, 6).Value = 52
$ws.Cells.Item(30, 6).Value = 1029
$ws.Cells.Item(33, 6).Value = 1011
$ws.Cells.Item(38, 6).Value = 2849
$ws.Cells.Item(41, 6).Value = 3376
$ws.Cells.Item(42, 6).Value = 1091
$ws.Cells.Item(45, 6).Value = 219
$ws.Cells.Item(49, 6).Value = 742
$ws.Cells.Item(51, 6).Value = 48

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 148367
$ws.Cells.Item(8, 6).Value = 148367
$ws.Cells.Item(15, 6).Value = 375
$ws.Cells.Item(26, 6).Value = 12

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(9, 6).Value = 704
$ws.Cells.Item(11, 6).Value = 2384

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 704
$ws.Cells.Item(7, 6).Value = 2384
$ws.Cells.Item(9, 6).Value = 533
$ws.Cells.Item(10, 6).Value = 148367
$ws.Cells.Item(12, 6).Value = 731
$ws.Cells.Item(17, 6).Value = 1764
$ws.Cells.Item(18, 6).Value = 1177
$ws.Cells.Item(20, 6).Value = 1697
$ws.Cells.Item(22, 6).Value = 25
$ws.Cells.Item(26, 6).Value = 1177
$ws.Cells.Item(27, 6).Value = 1783
$ws.Cells.Item(28, 6).Value = 1783
$ws.Cells.Item(29, 6).Value = 1185
$ws.Cells.Item(30, 6).Value = 165
$ws.Cells.Item(31, 6).Value = 1665
$ws.Cells.Item(32, 6).Value = 635
$ws.Cells.Item(33, 6).Value = 4497
$ws.Cells.Item(44, 6).Value = 2849
$ws.Cells.Item(47, 6).Value = 1091
$ws.Cells.Item(49, 6).Value = 219
$ws.Cells.Item(50, 6).Value = 703
